$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "ИТОГО"
$ws.Range("C7").Formula = "=SUM(C2:C6)"
$ws.Range("D7").Formula = "=SUM(D2:D6)"
$ws.Range("E7").Formula = "=SUM(E2:E6)"
